$d = $word.ActiveDocument

# --- Change 1: merge "Tests that a " + "nonnumerical" + " candidate size results in a "
#     into a single run with the combined text.
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "Tests that a nonnumerical candidate size results in a ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tests that a nonnumerical candidate size results in a ", 2)

# --- Change 2: merge " with a" + " nonnumerical" + " candidate size results in a "
#     into a single run with the combined text (keeps the run's rPr/formatting).
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    " with a nonnumerical candidate size results in a ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " with a nonnumerical candidate size results in a ", 2)

# --- Change 3: remove the whole "line: 2" paragraph (runs + paragraph mark).
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "line: 2*") {
        $p.Range.Delete()
    }
}
